$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '66.510.74' }
    @{ Cell = 'E2'; Value = '  -0.38%  ' }
    @{ Cell = 'D3'; Value = '3.517.99' }
    @{ Cell = 'E3'; Value = '  -3.30%  ' }
    @{ Cell = 'E4'; Value = '  +0.07%  ' }
    @{ Cell = 'D5'; Value = '606.21' }
    @{ Cell = 'E5'; Value = '  -0.72%  ' }
    @{ Cell = 'D6'; Value = '143.55' }
    @{ Cell = 'E6'; Value = '  -4.32%  ' }
    @{ Cell = 'D7'; Value = '3.517.75' }
    @{ Cell = 'E7'; Value = '  -3.33%  ' }
    @{ Cell = 'D8'; Value = '1.00' }
    @{ Cell = 'E8'; Value = '  -0.18%  ' }
    @{ Cell = 'D9'; Value = '0.507' }
    @{ Cell = 'E9'; Value = '  +3.44%  ' }
    @{ Cell = 'D10'; Value = '7.69' }
    @{ Cell = 'E10'; Value = '  -3.40%  ' }
    @{ Cell = 'E11'; Value = '  -5.27%  ' }
    @{ Cell = 'D12'; Value = '0.405' }
    @{ Cell = 'E12'; Value = '  -3.11%  ' }
    @{ Cell = 'D13'; Value = '4.111.54' }
    @{ Cell = 'E13'; Value = '  -3.29%  ' }
    @{ Cell = 'E14'; Value = '  -7.25%  ' }
    @{ Cell = 'D15'; Value = '28.66' }
    @{ Cell = 'E15'; Value = '  -4.64%  ' }
    @{ Cell = 'D16'; Value = '3.522.85' }
    @{ Cell = 'E16'; Value = '  -2.66%  ' }
    @{ Cell = 'E17'; Value = '  -0.11%  ' }
    @{ Cell = 'D18'; Value = '66.386.64' }
    @{ Cell = 'E18'; Value = '  -0.65%  ' }
    @{ Cell = 'D19'; Value = '10.74' }
    @{ Cell = 'E19'; Value = '  -7.52%  ' }
    @{ Cell = 'E20'; Value = '  -4.17%  ' }
    @{ Cell = 'D21'; Value = '14.59' }
    @{ Cell = 'E21'; Value = '  -4.05%  ' }
    @{ Cell = 'D22'; Value = '423.09' }
    @{ Cell = 'E22'; Value = '  -1.57%  ' }
    @{ Cell = 'D23'; Value = '0.589' }
    @{ Cell = 'E23'; Value = '  -5.52%  ' }
    @{ Cell = 'D24'; Value = '76.97' }
    @{ Cell = 'E24'; Value = '  -2.36%  ' }
    @{ Cell = 'D25'; Value = '3.667.07' }
    @{ Cell = 'E25'; Value = '  -2.90%  ' }
    @{ Cell = 'D26'; Value = '0.999' }
    @{ Cell = 'E26'; Value = '  -0.09%  ' }
    @{ Cell = 'D27'; Value = '0.0000113' }
    @{ Cell = 'E27'; Value = '  -8.02%  ' }
    @{ Cell = 'E28'; Value = '  -2.76%  ' }
    @{ Cell = 'D29'; Value = '7.81' }
    @{ Cell = 'E29'; Value = '  -7.07%  ' }
    @{ Cell = 'D30'; Value = '8.89' }
    @{ Cell = 'E30'; Value = '  -7.01%  ' }
    @{ Cell = 'E31'; Value = '  +0.05%  ' }
    @{ Cell = 'D32'; Value = '3.526.72' }
    @{ Cell = 'E32'; Value = '  -2.91%  ' }
    @{ Cell = 'D33'; Value = '0.153' }
    @{ Cell = 'E33'; Value = '  -3.47%  ' }
    @{ Cell = 'D34'; Value = '24.18' }
    @{ Cell = 'E34'; Value = '  -5.32%  ' }
    @{ Cell = 'E35'; Value = '  +0.01%  ' }
    @{ Cell = 'E36'; Value = '  -10.68%  ' }
    @{ Cell = 'D37'; Value = '7.53' }
    @{ Cell = 'E37'; Value = '  -5.00%  ' }
    @{ Cell = 'E38'; Value = '  -5.62%  ' }
    @{ Cell = 'D39'; Value = '173.35' }
    @{ Cell = 'E39'; Value = '  -2.06%  ' }
    @{ Cell = 'D40'; Value = '5.18' }
    @{ Cell = 'E40'; Value = '  -9.25%  ' }
    @{ Cell = 'D41'; Value = '0.0808' }
    @{ Cell = 'E41'; Value = '  -6.57%  ' }
    @{ Cell = 'E42'; Value = '  -5.60%  ' }
    @{ Cell = 'D43'; Value = '0.850' }
    @{ Cell = 'E43'; Value = '  -5.96%  ' }
    @{ Cell = 'E44'; Value = '  -0.87%  ' }
    @{ Cell = 'D45'; Value = '1.78' }
    @{ Cell = 'E45'; Value = '  -7.51%  ' }
    @{ Cell = 'E46'; Value = '  +0.11%  ' }
    @{ Cell = 'D47'; Value = '2.36' }
    @{ Cell = 'E47'; Value = '  -8.75%  ' }
    @{ Cell = 'D48'; Value = '7.05' }
    @{ Cell = 'E48'; Value = '  -2.43%  ' }
    @{ Cell = 'B49'; Value = 'EnergySwap' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' }
    @{ Cell = 'D49'; Value = '23.00' }
    @{ Cell = 'E49'; Value = '  -4.28%  ' }
    @{ Cell = 'B50'; Value = 'ONDO' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo' }
    @{ Cell = 'D50'; Value = '1.11' }
    @{ Cell = 'E50'; Value = '  -6.37%  ' }
    @{ Cell = 'D51'; Value = '0.901' }
    @{ Cell = 'E51'; Value = '  -6.48%  ' }
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    $c.NumberFormat = "@"
    $c.Value = $u.Value
    $c.Style = "Normal"
}